# Auto-generated edit script applying numeric corrections to profit calculation rows
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4495.2383
$ws.Range("I76").Value = 4812.5
$ws.Range("J76").Value = 3480
$ws.Range("K76").Value = 4812.5
$ws.Range("L76").Value = 3480
$ws.Range("M76").Value = -4497.5
$ws.Range("N76").Value = -4110

$ws.Range("H79").Value = 4495.2383
$ws.Range("I79").Value = 4812.5
$ws.Range("J79").Value = 3480
$ws.Range("K79").Value = 4812.5
$ws.Range("L79").Value = 3480
$ws.Range("M79").Value = -3720.5
$ws.Range("N79").Value = -5664

$ws.Range("H138").Value = 3098.0625
$ws.Range("I138").Value = 3181.5
$ws.Range("J138").Value = 3048
$ws.Range("K138").Value = 9544.5
$ws.Range("L138").Value = 9144
$ws.Range("M138").Value = -4404.5
$ws.Range("N138").Value = -19424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 13891621
$ws.Range("I61").Value = 23811776
$ws.Range("K61").Value = 23811776
$ws.Range("M61").Value = -23811564

$ws.Range("H74").Value = 960
$ws.Range("I74").Value = 847.6923
$ws.Range("K74").Value = 847.6923
$ws.Range("M74").Value = 26.30769999999995

$ws.Range("H77").Value = 960
$ws.Range("I77").Value = 847.6923
$ws.Range("K77").Value = 4238.4615
$ws.Range("M77").Value = 129.5384999999997

$ws.Range("H132").Value = 7511.826
$ws.Range("I132").Value = 12403
$ws.Range("J132").Value = 4903.2
$ws.Range("K132").Value = 37209
$ws.Range("L132").Value = 14709.6
$ws.Range("M132").Value = -34679
$ws.Range("N132").Value = -19769.6

$ws.Range("H136").Value = 13891621
$ws.Range("I136").Value = 23811776
$ws.Range("K136").Value = 71435328
$ws.Range("M136").Value = -71432778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2866.6667
$ws.Range("I105").Value = 2540
$ws.Range("J105").Value = 4500
$ws.Range("K105").Value = 2540
$ws.Range("L105").Value = 4500
$ws.Range("M105").Value = -793
$ws.Range("N105").Value = -7994

$ws.Range("H119").Value = 22111
$ws.Range("J119").Value = 22111
$ws.Range("L119").Value = 22111
$ws.Range("N119").Value = -31787

$ws.Range("H120").Value = 67880.5
$ws.Range("J120").Value = 67880.5
$ws.Range("L120").Value = 67880.5
$ws.Range("N120").Value = -77556.5

$ws.Range("H134").Value = 2251.6572
$ws.Range("I134").Value = 1976.9333
$ws.Range("K134").Value = 5930.7999
$ws.Range("M134").Value = -3395.7999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1349.2642
$ws.Range("I31").Value = 1072.5135
$ws.Range("J31").Value = 1989.25
$ws.Range("K31").Value = 1072.5135
$ws.Range("L31").Value = 1989.25
$ws.Range("M31").Value = -777.5135
$ws.Range("N31").Value = -2579.25

$ws.Range("H34").Value = 1349.2642
$ws.Range("I34").Value = 1072.5135
$ws.Range("J34").Value = 1989.25
$ws.Range("K34").Value = 1072.5135
$ws.Range("L34").Value = 1989.25
$ws.Range("M34").Value = -870.5135
$ws.Range("N34").Value = -2393.25

$ws.Range("H58").Value = 2478.5652
$ws.Range("I58").Value = 2229.3333
$ws.Range("J58").Value = 2750.4546
$ws.Range("K58").Value = 2229.3333
$ws.Range("L58").Value = 2750.4546
$ws.Range("M58").Value = -2026.3333
$ws.Range("N58").Value = -3156.4546

$ws.Range("H97").Value = 39800
$ws.Range("J97").Value = 39800
$ws.Range("L97").Value = 39800
$ws.Range("N97").Value = -41782

$ws.Range("H132").Value = 7938502.5
$ws.Range("I132").Value = 1463.75
$ws.Range("J132").Value = 18521220
$ws.Range("K132").Value = 4391.25
$ws.Range("L132").Value = 55563660
$ws.Range("M132").Value = -1861.25
$ws.Range("N132").Value = -55568720

$ws.Range("H134").Value = 1261.6842
$ws.Range("I134").Value = 813.2308
$ws.Range("K134").Value = 2439.6924
$ws.Range("M134").Value = 95.30760000000009

$ws.Range("H136").Value = 2478.5652
$ws.Range("I136").Value = 2229.3333
$ws.Range("J136").Value = 2750.4546
$ws.Range("K136").Value = 6687.999899999999
$ws.Range("L136").Value = 8251.363799999999
$ws.Range("M136").Value = -4137.999899999999
$ws.Range("N136").Value = -13351.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1011.9756
$ws.Range("I131").Value = 749.2308
$ws.Range("J131").Value = 1133.9642
$ws.Range("K131").Value = 2247.6924
$ws.Range("L131").Value = 3401.8926
$ws.Range("M131").Value = 2792.3076
$ws.Range("N131").Value = -13481.8926

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12579614
$ws.Range("I80").Value = 29281380
$ws.Range("J80").Value = 53290.25
$ws.Range("K80").Value = 29281380
$ws.Range("L80").Value = 53290.25
$ws.Range("M80").Value = -29280382
$ws.Range("N80").Value = -55286.25

$ws.Range("H83").Value = 12579614
$ws.Range("I83").Value = 29281380
$ws.Range("J83").Value = 53290.25
$ws.Range("K83").Value = 146406900
$ws.Range("L83").Value = 266451.25
$ws.Range("M83").Value = -146401908
$ws.Range("N83").Value = -276435.25

$ws.Range("H122").Value = 1668.9
$ws.Range("I122").Value = 1298.4286
$ws.Range("J122").Value = 2533.3333
$ws.Range("K122").Value = 3895.2858
$ws.Range("L122").Value = 7599.999899999999
$ws.Range("M122").Value = -1445.2858
$ws.Range("N122").Value = -12499.9999

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 3210.5715
$ws.Range("I132").Value = 2971.077
$ws.Range("J132").Value = 3599.75
$ws.Range("K132").Value = 8913.231
$ws.Range("L132").Value = 10799.25
$ws.Range("M132").Value = -6383.231
$ws.Range("N132").Value = -15859.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3825.628
$ws.Range("I132").Value = 3712.2632
$ws.Range("J132").Value = 3915.375
$ws.Range("K132").Value = 11136.7896
$ws.Range("L132").Value = 11746.125
$ws.Range("M132").Value = -8606.7896
$ws.Range("N132").Value = -16806.125

$ws.Range("H136").Value = 1128.8695
$ws.Range("I136").Value = 1152.909
$ws.Range("J136").Value = 600
$ws.Range("K136").Value = 3458.727
$ws.Range("L136").Value = 1800
$ws.Range("M136").Value = -908.7270000000003
$ws.Range("N136").Value = -6900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 33333.332
$ws.Range("J39").Value = 33333.332
$ws.Range("L39").Value = 33333.332
$ws.Range("N39").Value = -34159.332

$ws.Range("H46").Value = 59879
$ws.Range("J46").Value = 59879
$ws.Range("L46").Value = 59879
$ws.Range("N46").Value = -60341

$ws.Range("H134").Value = 59879
$ws.Range("J134").Value = 59879
$ws.Range("L134").Value = 179637
$ws.Range("N134").Value = -184707

$ws.Range("H136").Value = 2381.5
$ws.Range("I136").Value = 1991
$ws.Range("J136").Value = 3162.5
$ws.Range("K136").Value = 5973
$ws.Range("L136").Value = 9487.5
$ws.Range("M136").Value = -3423
$ws.Range("N136").Value = -14587.5
